# March 24 update 3
# Adds three new columns (M: renewd, N: PlanID, O: iteration) to the
# existing bldg sheet, filling every data row (2-13) with the new
# "before"/20150274/9 values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row entries, matching the existing header style (bold,
# bordered, centered) by copying the format from the last existing
# header cell (L1).
$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"

$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# New data values for every existing data row.
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 13).Value = "before"      # M - renewd
    $ws.Cells.Item($r, 14).Value = 20150274      # N - PlanID
    $ws.Cells.Item($r, 15).Value = 9             # O - iteration
}
